$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 142 (run_id 141)
$ws.Cells.Item(142, 1).Value = 141
$ws.Cells.Item(142, 2).Value = 1
$ws.Cells.Item(142, 3).Value = "2024-06-17 22:11:30"
$ws.Cells.Item(142, 4).Value = 200
$ws.Cells.Item(142, 5).Value = 14

# Row 143 (run_id 142)
$ws.Cells.Item(143, 1).Value = 142
$ws.Cells.Item(143, 2).Value = 2
$ws.Cells.Item(143, 3).Value = "2024-06-17 22:11:31"
$ws.Cells.Item(143, 4).Value = 200
$ws.Cells.Item(143, 5).Value = 2
